$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row suffixes: "_old" -> "_FV2304", "_new" -> "_FV2310" ---
# Columns A-J ("*_old" headers, formerly the FV2304 input file's columns)
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
# Columns L-U ("*_new" headers, formerly the FV2310 input file's columns)
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $fv2304Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2304Headers[$i]
}

# column 11 ("diff") stays untouched

for ($i = 0; $i -lt $fv2310Headers.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $fv2310Headers[$i]
}

# --- 2. Turn the used range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U66"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split after row 1) ---
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
